# Generate Report for Handback
# - Mark the f8b0c4c0 row's zh-cn/de-de handback status as failed, and
#   record the detailed error message in the "Error Detail" column of the
#   per-language report sheets. Also widen that column so the message is
#   readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text for the f8b0c4c0-... file (row 3) changes from
# "Ready for handoff" to "Handback transform failed" everywhere it is
# shown: the Overview sheet (zh-cn/de-de status columns) and each
# per-language sheet's own Status column.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Widen the "Error Detail" column (column P) on both language sheets so
# the new, longer error message is readable. ColumnWidth is specified in
# characters and round-trips through a pixel-based rounding step, so the
# value below is chosen to land on an on-disk column width of exactly 40.
$targetColumnWidth = 40 - (5 / 6)
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColumnWidth

# Record the handback/handoff filename mismatch error detail for the
# f8b0c4c0-... row (row 3) on each language sheet.
$wsZhCn.Range("P3").Value = "Handback file name: xzvpxxp5.phx is different with handoff file name: f8b0c4c0-fb3c-4d67-9e33-8d48bad0818a.c73a4cc6363e3c8332821ef270ffa34dd8647f28.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: xzvpxxp5.phx is different with handoff file name: f8b0c4c0-fb3c-4d67-9e33-8d48bad0818a.c73a4cc6363e3c8332821ef270ffa34dd8647f28.de-de."
